# Update column F ("dSF") values for specific rows to reflect a data re-pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    7  = -6
    21 = -5
    23 = -3
    24 = -1
    28 = 8
    29 = -4
    30 = 8
    33 = 1
    38 = -2
    39 = -4
    40 = -5
    42 = 2
    43 = -7
    46 = -2
    52 = 1
    58 = 0
    59 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
